$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("N1_D40")
$ws.Range("E2").Value = 0.112
$ws.Range("F2").Value = 17.04
$ws.Range("E3").Value = 0.055
$ws.Range("F3").Value = 17.18
$ws.Range("E4").Value = 0.055
$ws.Range("F4").Value = 17.18
$ws.Range("E5").Value = 0.052
$ws.Range("F5").Value = 16.8
$ws.Range("F6").Value = 17.29
$ws.Range("E7").Value = 0.049
$ws.Range("F7").Value = 16.95
$ws.Range("E8").Value = 0.053
$ws.Range("F8").Value = 17.12
$ws.Range("E9").Value = 0.054
$ws.Range("F9").Value = 17.07
$ws.Range("E10").Value = 0.053
$ws.Range("F10").Value = 17.01
$ws.Range("E11").Value = 0.055
$ws.Range("F11").Value = 17.11
$ws.Range("E12").Value = 0.05920000000000001
$ws.Range("F12").Value = 17.075

$ws = $wb.Worksheets.Item("N1_D60")
$ws.Range("E2").Value = 0.08500000000000001
$ws.Range("F2").Value = 19.27
$ws.Range("E3").Value = 0.083
$ws.Range("F3").Value = 18.89
$ws.Range("E4").Value = 0.083
$ws.Range("F4").Value = 18.86
$ws.Range("E5").Value = 0.083
$ws.Range("F5").Value = 19.19
$ws.Range("E6").Value = 0.08400000000000001
$ws.Range("F6").Value = 18.94
$ws.Range("E7").Value = 0.081
$ws.Range("F7").Value = 19.14
$ws.Range("E8").Value = 0.079
$ws.Range("F8").Value = 18.97
$ws.Range("E9").Value = 0.083
$ws.Range("F9").Value = 18.7
$ws.Range("E10").Value = 0.08500000000000001
$ws.Range("F10").Value = 19.25
$ws.Range("E11").Value = 0.08500000000000001
$ws.Range("F11").Value = 19.77
$ws.Range("E12").Value = 0.08309999999999999
$ws.Range("F12").Value = 19.098

$ws = $wb.Worksheets.Item("N1_D80")
$ws.Range("F2").Value = 26.01
$ws.Range("E3").Value = 0.14
$ws.Range("F3").Value = 25.95
$ws.Range("E4").Value = 0.14
$ws.Range("F4").Value = 26.02
$ws.Range("E5").Value = 0.141
$ws.Range("F5").Value = 26.04
$ws.Range("E6").Value = 0.139
$ws.Range("F6").Value = 26.06
$ws.Range("E7").Value = 0.139
$ws.Range("F7").Value = 25.91
$ws.Range("E8").Value = 0.141
$ws.Range("F8").Value = 26
$ws.Range("E9").Value = 0.141
$ws.Range("F9").Value = 26.08
$ws.Range("F10").Value = 25.96
$ws.Range("E11").Value = 0.138
$ws.Range("F11").Value = 25.98
$ws.Range("E12").Value = 0.1396
$ws.Range("F12").Value = 26.001

$ws = $wb.Worksheets.Item("N1_D100")
$ws.Range("E2").Value = 0.148
$ws.Range("F2").Value = 29.62
$ws.Range("E3").Value = 0.144
$ws.Range("F3").Value = 29.57
$ws.Range("E4").Value = 0.15
$ws.Range("F4").Value = 31.18
$ws.Range("E5").Value = 0.146
$ws.Range("F5").Value = 29.51
$ws.Range("E6").Value = 0.148
$ws.Range("F6").Value = 29.72
$ws.Range("E7").Value = 0.148
$ws.Range("F7").Value = 31.43
$ws.Range("E8").Value = 0.149
$ws.Range("F8").Value = 29.77
$ws.Range("E9").Value = 0.149
$ws.Range("F9").Value = 29.5
$ws.Range("E10").Value = 0.147
$ws.Range("F10").Value = 29.34
$ws.Range("E11").Value = 0.149
$ws.Range("F11").Value = 29.56
$ws.Range("E12").Value = 0.1478
$ws.Range("F12").Value = 29.92
